$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student rows (6-10), mirroring the style of the last existing data
# row (row 5) for column A (style index 3 in the original sheet).
$ws.Range("A5:I5").Copy() | Out-Null
$ws.Range("A6:I10").PasteSpecial(-4122) | Out-Null

$rows = @(
    @{ Nim = 9999999994; Nama = "das";   C = 90; D = 100; E = 88;  F = 100; G = 50;  H = 50;  I = 88 },
    @{ Nim = 9999999995; Nama = "eq";    C = 50; D = 88;  E = 50;  F = 90;  G = 88;  H = 88;  I = 88 },
    @{ Nim = 9999999996; Nama = "fsda";  C = 50; D = 90;  E = 90;  F = 100; G = 100; H = 90;  I = 100 },
    @{ Nim = 9999999997; Nama = "dsd";   C = 90; D = 100; E = 100; F = 88;  G = 88;  H = 100; I = 90 },
    @{ Nim = 9999999998; Nama = "sdewq"; C = 90; D = 100; E = 100; F = 88;  G = 88;  H = 100; I = 90 }
)

$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Nim
    $ws.Cells.Item($r, 2).Value = $row.Nama
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $r = $r + 1
}

# Drop the stray outline-level-row metadata that no longer reflects any
# real row grouping, by normalising every row's OutlineLevel to 0.
for ($i = 1; $i -le 10; $i++) {
    $ws.Rows.Item($i).OutlineLevel = 0
}

# Update the selection to match the post-edit UI state.
$ws.Range("C10:I10").Select() | Out-Null
